$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 2020 -> 2021 model-year refresh, with a few MSRP (column D) updates ---

# Row 29: ES 350 Luxury-ish trim — year only
$ws.Range("C29").Value = 2021

# Row 30: year + MSRP bump
$ws.Range("C30").Value = 2021
$ws.Range("D30").Value = 56190

# Row 31: year only
$ws.Range("C31").Value = 2021

# Row 45: year + MSRP bump
$ws.Range("C45").Value = 2021
$ws.Range("D45").Value = 76000

# Row 46: year + MSRP bump
$ws.Range("C46").Value = 2021
$ws.Range("D46").Value = 79250

# Row 47: year + MSRP drop
$ws.Range("C47").Value = 2021
$ws.Range("D47").Value = 79600

# Row 48: year + MSRP drop
$ws.Range("C48").Value = 2021
$ws.Range("D48").Value = 82850

# --- New row 99: UX 250h AWD Black Line Special Edition (deeplink/msrp test row) ---
$ws.Range("B99").Value = "UX 250h AWD BLACK LINE SPECIAL EDITION"
$ws.Range("C99").Value = 2021
$ws.Range("D99").Value = "COMING SOON"

# Match formatting of neighboring rows: D99 picks up the plain-number style
# used by D73:D78, E99 picks up the same style already used on E98.
$ws.Range("D73").Copy()
$ws.Range("D99").PasteSpecial(-4122)
$ws.Range("D99").Value = "COMING SOON"

$ws.Range("E98").Copy()
$ws.Range("E99").PasteSpecial(-4122)
$ws.Range("E99").Value = 1025

# --- View state: scroll position + active selection moved up the sheet ---
$ws.Range("D49").Select()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
